$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The old row 8 (header-like row with Index/Date/Value labels) is removed;
# data rows shift up by one (old row 9 becomes new row 8, etc.)
$ws.Rows("8:8").Delete()

# Column C now holds a running-count formula instead of constant 0 values.
# C8 (old C9) starts the sequence at -1, C9 (old C10) adds 1 to the previous
# cell, and the remaining cells (C10:C17) share that same relative formula.
$ws.Range("C8").Formula = "=-1"
$ws.Range("C9").Formula = "=C8+1"
$ws.Range("C10:C17").Formula = "=C9+1"

# Update the active selection to match the new state of the sheet.
$ws.Range("E10").Select()
